$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cid values in column A for rows 7-14 per Tube Data separation
$ws.Range("A7").Value = 7101
$ws.Range("A8").Value = 7102
$ws.Range("A9").Value = 7103
$ws.Range("A10").Value = 7104
$ws.Range("A11").Value = 7105
$ws.Range("A12").Value = 7106
$ws.Range("A13").Value = 7107
$ws.Range("A14").Value = 7108

# Update the active selection on the sheet
$ws.Range("B6").Select()
